# Fill in the results for the 13ra Giornata lega / 14ra Giornata serie a
# round (rows 30-33 on the "Calendario" sheet): these four matches
# previously showed placeholder 0-0 scores and "-" results; the
# commit uploads the real final scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendario")

# Match 1 (row 30): Scroto FC vs Canile Comunale Di Merate
$ws.Range("B30").Value = 69.5
$ws.Range("C30").Value = 78
$ws.Range("E30").Value = "1-4"

# Match 2 (row 31): FC Tumori vs Chianti e Pianti
$ws.Range("B31").Value = 72.5
$ws.Range("C31").Value = 67.5
$ws.Range("E31").Value = "2-1"

# Match 3 (row 32): Black Gay United vs FC ETTANERA
$ws.Range("B32").Value = 78.5
$ws.Range("C32").Value = 72
$ws.Range("E32").Value = "4-2"

# Match 4 (row 33): CSKA LA RISSA vs ReAlcolizzati
$ws.Range("B33").Value = 67.5
$ws.Range("C33").Value = 70
$ws.Range("E33").Value = "1-2"
